# Update "想去人数" (want-to-go counts) figures on the 展览 and 全部类型 sheets,
# reflecting a refreshed data pull (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 72
$wsExhibit.Range("F11").Value = 10105
$wsExhibit.Range("F13").Value = 262
$wsExhibit.Range("F15").Value = 628
$wsExhibit.Range("F17").Value = 12143

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 72
$wsAll.Range("F12").Value = 10105
$wsAll.Range("F14").Value = 262
$wsAll.Range("F16").Value = 628
$wsAll.Range("F18").Value = 12143
